# Weekly Understanding Portfolio 5 / JSON
#
# Inserts a large block of new weekly-log paragraphs after the
# "1:00 I asked Henry what JSON was" paragraph, moves the "Portfolio code "
# paragraph to right after that point (duplicating its old text earlier in
# the doc), and relocates the _GoBack bookmark into its own empty paragraph
# at the very end of the inserted block (removing it from its old spot).

$d = $word.ActiveDocument

# ---- helpers -------------------------------------------------------------

# Insert a brand-new, *truly* empty paragraph right after paragraph
# $paraIndex and return its (1-based) paragraph index. We "touch" the new
# paragraph with throwaway text and then clear it; that collapses away the
# phantom run that Word would otherwise leave behind for the paragraph
# mark, so later InsertAfter/InsertBreak/Hyperlinks.Add calls land cleanly
# with no stray empty <w:r/>.
function New-EmptyParagraphAfter {
    param($paraIndex)
    $doc = $word.ActiveDocument
    $doc.Paragraphs($paraIndex).Range.InsertParagraphAfter()
    $newIndex = $paraIndex + 1
    $touch = $doc.Paragraphs($newIndex).Range
    $touch.Text = "x"
    $afterTouch = $doc.Paragraphs($newIndex).Range
    $clear = $doc.Range($afterTouch.Start, $afterTouch.End - 1)
    $clear.Text = ""
    return $newIndex
}

# Append a new run of plain text to the end of paragraph $paraIndex
# (just before its paragraph mark).
function Add-RunText {
    param($paraIndex, $text)
    $doc = $word.ActiveDocument
    $para = $doc.Paragraphs($paraIndex).Range
    $tail = $doc.Range($para.End - 1, $para.End - 1)
    $tail.InsertAfter($text)
}

# Append a line-break run to the end of paragraph $paraIndex.
function Add-RunBreak {
    param($paraIndex)
    $doc = $word.ActiveDocument
    $para = $doc.Paragraphs($paraIndex).Range
    $tail = $doc.Range($para.End - 1, $para.End - 1)
    $tail.InsertBreak(6) | Out-Null
}

# Turn (empty) paragraph $paraIndex into a hyperlink paragraph.
function Add-ParagraphHyperlink {
    param($paraIndex, $url)
    $doc = $word.ActiveDocument
    $para = $doc.Paragraphs($paraIndex).Range
    $doc.Hyperlinks.Add($para, $url) | Out-Null
}

# Drop a bookmark named $name, collapsed, inside (empty) paragraph
# $paraIndex.
function Add-ParagraphBookmark {
    param($paraIndex, $name)
    $doc = $word.ActiveDocument
    $para = $doc.Paragraphs($paraIndex).Range
    $doc.Bookmarks.Add($name, $para) | Out-Null
}

# ---- locate the anchor paragraph -----------------------------------------

# "1:00 I asked Henry what JSON was" is paragraph 4 in the original
# document; find it by content instead of hard-coding the index so the
# script is resilient to minor structural differences.
$anchor = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "1:00 I asked Henry what JSON was*") {
        $anchor = $i
    }
}

# The _GoBack bookmark currently sits at the end of the anchor paragraph;
# remove it from there (it will be re-added in its own paragraph at the
# end of the newly inserted block).
$goBack = $d.Bookmarks
if ($goBack.Exists("_GoBack")) {
    $goBack.Item("_GoBack").Delete()
}

$idx = $anchor

# ---- "Portfolio code " (duplicated earlier in the doc) -------------------
$idx = New-EmptyParagraphAfter $idx
Add-RunText $idx "Portfolio code "

# ---- blank paragraph -------------------------------------------------------
$idx = New-EmptyParagraphAfter $idx

# ---- WEEK 2/7/2015 ---------------------------------------------------------
$idx = New-EmptyParagraphAfter $idx
Add-RunText $idx "WEEK 2/7/2015"

# ---- At around 6 sec I started teaching -----------------------------------
$idx = New-EmptyParagraphAfter $idx
Add-RunText $idx "At around 6 sec I started teaching"

# ---- No questions for me ---------------------------------------------------
$idx = New-EmptyParagraphAfter $idx
Add-RunText $idx "No questions for me"

# ---- No answers <br/> Did javascript refresher. ---------------------------
$idx = New-EmptyParagraphAfter $idx
Add-RunText $idx "No answers "
Add-RunBreak $idx
Add-RunText $idx "Did"
Add-RunText $idx " "
Add-RunText $idx "javascript"
Add-RunText $idx " refresher."

# ---- Portfolio code is on github ------------------------------------------
$idx = New-EmptyParagraphAfter $idx
Add-RunText $idx "Portfolio code is on "
Add-RunText $idx "github"

# ---- Ch7_20.html ------------------------------------------------------------
$idx = New-EmptyParagraphAfter $idx
Add-RunText $idx "Ch7_20.html"

# ---- Palindrome javascript example... --------------------------------------
$idx = New-EmptyParagraphAfter $idx
Add-RunText $idx "Palindrome "
Add-RunText $idx "javascript"
Add-RunText $idx " example with output in the format of a table"

# ---- MyMenu3.html -----------------------------------------------------------
$idx = New-EmptyParagraphAfter $idx
Add-RunText $idx "MyMenu3.html"

# ---- Now shows a date on it via javascript getYear function -----------------
$idx = New-EmptyParagraphAfter $idx
Add-RunText $idx "Now shows a date on it via "
Add-RunText $idx "javascript"
Add-RunText $idx " "
Add-RunText $idx "getYear"
Add-RunText $idx " function"

# ---- <br/>LESSON 5 -----------------------------------------------------------
$idx = New-EmptyParagraphAfter $idx
Add-RunBreak $idx
Add-RunText $idx "LESSON 5"

# ---- WEEK 02/14/2015 ---------------------------------------------------------
$idx = New-EmptyParagraphAfter $idx
Add-RunText $idx "WEEK 02/14/2015"

# ---- hyperlink: myjson.html --------------------------------------------------
$idx = New-EmptyParagraphAfter $idx
Add-ParagraphHyperlink $idx "http://shaferprojectsite.webatu.com/myjson.html"

# ---- hyperlink: W3SchoolsJSON.html -------------------------------------------
$idx = New-EmptyParagraphAfter $idx
Add-ParagraphHyperlink $idx "http://shaferprojectsite.webatu.com/W3SchoolsJSON.html"

# ---- final paragraph holding the relocated _GoBack bookmark -----------------
$idx = New-EmptyParagraphAfter $idx
Add-ParagraphBookmark $idx "_GoBack"

Write-Output ("Inserted through paragraph " + $idx + " of " + $d.Paragraphs.Count)
